$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of data (dormroom.jpg) in row 31
$ws.Range("A31").Value = "dormroom.jpg"
$ws.Range("B31").Value = "Picture of Male's Dormroom"
$ws.Range("C31").Value = "1280x1080px"
$ws.Range("D31").Value = "Original"
$ws.Range("E31").Value = "Placeholder"

# Update the view: scroll so row 7 is at top, and select E31
$ws.Range("E31").Select()
$excel.ActiveWindow.ScrollRow = 7
